# fall 24 week 10 inputs and lineup message improvements
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.26
$ws.Range("F2").Value = 0.83

$ws.Range("B3").Value = 1.59
$ws.Range("D3").Value = 1.38

$ws.Range("C4").Value = 1.47
$ws.Range("D4").Value = 1.32
$ws.Range("F4").Value = 1.09

$ws.Range("D5").Value = 1.34

$ws.Range("B6").Value = 2
$ws.Range("D6").Value = 1.5
$ws.Range("G6").Value = 0.95

$ws.Range("F7").Value = 1.48
$ws.Range("G7").Value = 1.15
